$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 205.875
$ws.Range("J31").Value = 107.833336
$ws.Range("L31").Value = 323.500008
$ws.Range("N31").Value = -783.500008

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 147091.56
$ws.Range("I132").Value = 169017.03
$ws.Range("J132").Value = 921.7778
$ws.Range("K132").Value = 507051.09
$ws.Range("L132").Value = 2765.3334
$ws.Range("M132").Value = -504521.09
$ws.Range("N132").Value = -7825.3334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 11779237
$ws.Range("I135").Value = 444.25
$ws.Range("J135").Value = 21202272
$ws.Range("K135").Value = 3998.25
$ws.Range("L135").Value = 190820448
$ws.Range("M135").Value = -1463.25
$ws.Range("N135").Value = -190825518

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 29167900
$ws.Range("I137").Value = 6667792.5
$ws.Range("J137").Value = 66668080
$ws.Range("K137").Value = 20003377.5
$ws.Range("L137").Value = 200004240
$ws.Range("M137").Value = -20000827.5
$ws.Range("N137").Value = -200009340

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2489.6382
$ws.Range("I138").Value = 2149.6428
$ws.Range("J138").Value = 2990.6843
$ws.Range("K138").Value = 6448.928400000001
$ws.Range("L138").Value = 8972.052899999999
$ws.Range("M138").Value = -1308.928400000001
$ws.Range("N138").Value = -19252.0529

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1339.5714
$ws.Range("I141").Value = 578.2222
$ws.Range("J141").Value = 2710
$ws.Range("K141").Value = 1734.6666
$ws.Range("L141").Value = 8130
$ws.Range("M141").Value = 3445.3334
$ws.Range("N141").Value = -18490

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4379.672
$ws.Range("I32").Value = 4355.121
$ws.Range("K32").Value = 4355.121
$ws.Range("M32").Value = -4068.121

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 63522.875
$ws.Range("I45").Value = 250693.25
$ws.Range("J45").Value = 1132.75
$ws.Range("K45").Value = 250693.25
$ws.Range("L45").Value = 1132.75
$ws.Range("M45").Value = -250316.25
$ws.Range("N45").Value = -1886.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5004472
$ws.Range("I74").Value = 7143650
$ws.Range("J74").Value = 13056.917
$ws.Range("K74").Value = 7143650
$ws.Range("L74").Value = 13056.917
$ws.Range("M74").Value = -7142776
$ws.Range("N74").Value = -14804.917

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 5004472
$ws.Range("I77").Value = 7143650
$ws.Range("J77").Value = 13056.917
$ws.Range("K77").Value = 35718250
$ws.Range("L77").Value = 65284.585
$ws.Range("M77").Value = -35713882
$ws.Range("N77").Value = -74020.58499999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1127.9375
$ws.Range("I110").Value = 953.5
$ws.Range("J110").Value = 1418.6666
$ws.Range("K110").Value = 953.5
$ws.Range("L110").Value = 1418.6666
$ws.Range("M110").Value = 1091.5
$ws.Range("N110").Value = -5508.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 935697.5
$ws.Range("I132").Value = 1114357.8
$ws.Range("J132").Value = 144487.72
$ws.Range("K132").Value = 3343073.4
$ws.Range("L132").Value = 433463.16
$ws.Range("M132").Value = -3340543.4
$ws.Range("N132").Value = -438523.16

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1456.3214
$ws.Range("I20").Value = 1432.4445
$ws.Range("J20").Value = 1499.3
$ws.Range("K20").Value = 1432.4445
$ws.Range("L20").Value = 1499.3
$ws.Range("M20").Value = -1185.4445
$ws.Range("N20").Value = -1993.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1670.2222
$ws.Range("I86").Value = 1452.9592
$ws.Range("J86").Value = 2133.087
$ws.Range("K86").Value = 1452.9592
$ws.Range("L86").Value = 2133.087
$ws.Range("M86").Value = -329.9592
$ws.Range("N86").Value = -4379.087

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1670.2222
$ws.Range("I89").Value = 1452.9592
$ws.Range("J89").Value = 2133.087
$ws.Range("K89").Value = 7264.796
$ws.Range("L89").Value = 10665.435
$ws.Range("M89").Value = -1648.796
$ws.Range("N89").Value = -21897.435

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1525
$ws.Range("I16").Value = 2125
$ws.Range("J16").Value = 925
$ws.Range("K16").Value = 2125
$ws.Range("L16").Value = 925
$ws.Range("M16").Value = -1838
$ws.Range("N16").Value = -1499

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 2778225.8
$ws.Range("I107").Value = 4629920.5
$ws.Range("J107").Value = 683.6667
$ws.Range("K107").Value = 4629920.5
$ws.Range("L107").Value = 683.6667
$ws.Range("M107").Value = -4628000.5
$ws.Range("N107").Value = -4523.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1525
$ws.Range("I113").Value = 2125
$ws.Range("J113").Value = 925
$ws.Range("K113").Value = 2125
$ws.Range("L113").Value = 925
$ws.Range("M113").Value = 45
$ws.Range("N113").Value = -5265

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1547.2273
$ws.Range("I132").Value = 1596.8572
$ws.Range("J132").Value = 505
$ws.Range("K132").Value = 4790.571599999999
$ws.Range("L132").Value = 1515
$ws.Range("M132").Value = -2260.571599999999
$ws.Range("N132").Value = -6575

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 41738.957
$ws.Range("I12").Value = 72.90909000000001
$ws.Range("J12").Value = 76994.84
$ws.Range("K12").Value = 218.72727
$ws.Range("L12").Value = 230984.52
$ws.Range("M12").Value = -45.72727000000003
$ws.Range("N12").Value = -231330.52

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 2000
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 2000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 6000
$ws.Range("N20").Value = -6454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3435.342
$ws.Range("I131").Value = 7562.857
$ws.Range("J131").Value = 2503.3225
$ws.Range("K131").Value = 22688.571
$ws.Range("L131").Value = 7509.967500000001
$ws.Range("M131").Value = -17648.571
$ws.Range("N131").Value = -17589.9675

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 3000
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1529.9667
$ws.Range("I102").Value = 1407.1666
$ws.Range("J102").Value = 1714.1666
$ws.Range("K102").Value = 1407.1666
$ws.Range("L102").Value = 1714.1666
$ws.Range("M102").Value = 214.8334
$ws.Range("N102").Value = -4958.1666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 554.4761999999999
$ws.Range("I107").Value = 377.83334
$ws.Range("J107").Value = 790
$ws.Range("K107").Value = 377.83334
$ws.Range("L107").Value = 790
$ws.Range("M107").Value = 1542.16666
$ws.Range("N107").Value = -4630

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3936.3635
$ws.Range("I126").Value = 3500
$ws.Range("J126").Value = 3980
$ws.Range("K126").Value = 10500
$ws.Range("L126").Value = 11940
$ws.Range("M126").Value = -8030
$ws.Range("N126").Value = -16880

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 12515.5
$ws.Range("I35").Value = 12515.5
$ws.Range("K35").Value = 12515.5
$ws.Range("M35").Value = -12179.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 1730.8334
$ws.Range("I122").Value = 1807.5
$ws.Range("J122").Value = 1577.5
$ws.Range("K122").Value = 5422.5
$ws.Range("L122").Value = 4732.5
$ws.Range("M122").Value = -2972.5
$ws.Range("N122").Value = -9632.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4824.567
$ws.Range("I136").Value = 4714.28
$ws.Range("K136").Value = 14142.84
$ws.Range("M136").Value = -11592.84

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5041.6284
$ws.Range("I132").Value = 6021.448
$ws.Range("J132").Value = 305.83334
$ws.Range("K132").Value = 18064.344
$ws.Range("L132").Value = 917.5000200000001
$ws.Range("M132").Value = -15534.344
$ws.Range("N132").Value = -5977.50002

# Remove cells that no longer exist in the target state
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M20").ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M4").ClearContents()
